$d = $word.ActiveDocument

# 1. Change "Temperature Measurement" to "Temp Service" -- temporarily keep a
#    trailing marker so the insertion point used below is not the very last
#    character of the paragraph (avoids an edge-case in point-bookmark
#    placement at a paragraph's final character position).
$d.Content.Find.Execute("Temperature Measurement", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Temp ServiceXXMARKERXX", 2)

# 2. Locate the end of "Temp Service" (before the temporary marker) and drop
#    the _GoBack bookmark there. Re-adding a bookmark with an existing name
#    relocates it, so this also removes it from its old spot near the
#    "More to be added" paragraph.
$found = $d.Content
$found.Find.Execute("Temp Service", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0)
$found.Collapse(0)
$d.Bookmarks.Add("_GoBack", $found) | Out-Null

# 3. Remove the temporary marker text.
$d.Content.Find.Execute("XXMARKERXX", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)

# 4. Merge "More to be added as we are still deciding as a group to split the
#    work" and the following standalone space run into a single run whose
#    text already carries the trailing space.
$d.Content.Find.Execute("More to be added as we are still deciding as a group to split the work ",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false, "More to be added as we are still deciding as a group to split the work ", 2)
